$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 23, pushing existing rows 23:119 down to 24:120.
# Excel's Insert() copies formatting from the row above (row 22), matching
# the date-format style (s="2") already present on column D.
$ws.Rows("23:23").Insert()

# Populate the newly inserted row 23 with the new data record.
$ws.Range("A23").Value = 10
$ws.Range("B23").Value = "Vega Modelo de Temuco"
$ws.Range("C23").Value = "La Araucanía"
$ws.Range("D23").Value = "2023-11-20"
$ws.Range("E23").Value = 9
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100107
$ws.Range("H23").Value = "Otros"
$ws.Range("I23").Value = 100107011
$ws.Range("J23").Value = "Tuna"
$ws.Range("K23").Value = "Sin especificar"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 110
$ws.Range("N23").Value = 40000
$ws.Range("O23").Value = 40000
$ws.Range("P23").Value = 40000
$ws.Range("Q23").Value = "$/caja 16 kilos"
$ws.Range("R23").Value = "Provincia de Los Andes"
$ws.Range("S23").Value = 2500
$ws.Range("T23").Value = 16
